# (feat) WIP: add single chart to slide
#
# The chart graphicFrame that was dropped onto the slide was still
# carrying PowerPoint's generic auto-generated name ("Diagramm 5" /
# "Chart 5"). Give it a proper, descriptive name now that it is the
# single chart on the slide.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shape = $s.Shapes.Item($i)
    if ($shape.HasChart) {
        $shape.Name = "StackedBars"
    }
}
